$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price and volume(1h) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.903.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4584'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3693'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07172'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8770'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07838'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.62'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.818.13'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.45%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.393'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.24'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.009'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008711'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.933.83'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.97'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.17'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.969'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.91%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.922'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08800'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.053'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7538'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.483'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.132'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.565'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.085'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01934'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05134'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.906'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.945'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4977'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1595'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.300'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4680'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.007'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.17'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.12'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.614'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06117'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.36'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.07%  '
